$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.897.56"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("E2").ClearFormats()

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.629.62"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E3").ClearFormats()

# Row 4 - TetherUSD
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E4").ClearFormats()

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.78"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("E5").ClearFormats()

# Row 6 - XRP
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("E6").ClearFormats()

# Row 7 - USDC
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E7").ClearFormats()

# Row 8 - Solana
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.33%  "
$ws.Range("E8").ClearFormats()

# Row 9 - Cardano
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.97%  "
$ws.Range("E9").ClearFormats()

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0880"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("E11").ClearFormats()

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.861.50"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("E12").ClearFormats()

# Row 13 - WrappedEther
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.623.72"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("E13").ClearFormats()

# Row 14 - Polkadot
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("E14").ClearFormats()

# Row 15 - Polygon
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.27%  "
$ws.Range("E15").ClearFormats()

# Row 16 - Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.56"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("E16").ClearFormats()

# Row 17 - WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.895.49"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("E17").ClearFormats()

# Row 18 - BitcoinCash
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.64"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("E18").ClearFormats()

# Row 19 - Chainlink
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.67"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("E19").ClearFormats()

# Row 20 - ShibaInu
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("E20").ClearFormats()

# Row 21 - Dai
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("E21").ClearFormats()

# Row 22 - Uniswap
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("E22").ClearFormats()

# Row 23 - Avalanche
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.22"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -5.44%  "
$ws.Range("E23").ClearFormats()

# Row 24 - Toncoin
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("E24").ClearFormats()

# Row 25 - Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.86"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.12%  "
$ws.Range("E25").ClearFormats()

# Row 26 - Cosmos
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("E26").ClearFormats()

# Row 27 - Stellar
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E27").ClearFormats()

# Row 28 - EthereumClassic
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("E28").ClearFormats()

# Row 29 - BinanceUSD
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E29").ClearFormats()

# Row 30 - PancakeSwap
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.19"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("E30").ClearFormats()

# Row 31 - Hedera
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("E31").ClearFormats()

# Row 32 - Filecoin
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("E32").ClearFormats()

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("E33").ClearFormats()

# Row 34 - Maker
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.400.98"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("E34").ClearFormats()

# Row 35 - LidoDAOToken
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("E35").ClearFormats()

# Row 36 - TrustWalletToken
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +8.25%  "
$ws.Range("E36").ClearFormats()

# Row 37 - HuobiToken
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("E37").ClearFormats()

# Row 38 - VeChain
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.48%  "
$ws.Range("E38").ClearFormats()

# Row 39 - ImmutableX
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.556"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("E39").ClearFormats()

# Row 40 - ARBITRUM
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.865"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.02%  "
$ws.Range("E40").ClearFormats()

# Row 41 - WEMIXToken
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E41").ClearFormats()

# Row 42 - PaxDollar
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("E42").ClearFormats()

# Row 43 - RenderToken
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("E43").ClearFormats()

# Row 44 - FraxShare
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("E44").ClearFormats()

# Row 45 - Aave
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.09"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("E45").ClearFormats()

# Row 46 - MXToken
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("E46").ClearFormats()

# Row 47 - RocketPoolETH
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.770.80"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("E47").ClearFormats()

# Row 48 - Quant
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.22"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("E48").ClearFormats()

# Row 49 - Algorand -> BabyDogeCoin (row shift)
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0104"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.71%  "
$ws.Range("E49").ClearFormats()

# Row 50 - Cronos -> Algorand (row shift)
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.102"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("E50").ClearFormats()

# Row 51 - EnergySwap -> Cronos (row shift)
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0504"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.42%  "
$ws.Range("E51").ClearFormats()
